$wb = $excel.ActiveWorkbook

# Add the new "PatientDetails" worksheet after the last existing sheet ("Attribute")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "PatientDetails"

# Cell values, entered in the same order the original schema-builder tool
# emitted them (controls shared-string table insertion order).
$ws.Range("D1").Value = "default"
$ws.Range("B3").Value = "InsuranceNumber"
$ws.Range("D3").Value = "123456789ABC"
$ws.Range("B4").Value = "DateOfBirth"
$ws.Range("C4").Value = "date"
$ws.Range("H2").Value = "'This is the Schema for Basic Tutorial"
$ws.Range("B5").Value = "Gender"
$ws.Range("B6").Value = "Weight"
$ws.Range("C6").Value = "decimal"
$ws.Range("H4").Value = "DateOfBirth docs"
$ws.Range("E1").Value = "values"
$ws.Range("E5").Value = "male, female"
$ws.Range("F1").Value = "unit.values"
$ws.Range("G1").Value = "unit.default"
$ws.Range("F6").Value = "g, kg"
$ws.Range("G6").Value = "kg"
$ws.Range("B2").Value = "PatientDetails"

# Remaining header / row values that reuse shared strings already present
# in the workbook (no new shared-string entries required for these).
$ws.Range("A1").Value = "class"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "type"
$ws.Range("H1").Value = "documentation"
$ws.Range("A2").Value = "struct"
$ws.Range("A3").Value = "attribute"
$ws.Range("C3").Value = "string"
$ws.Range("A4").Value = "field"
$ws.Range("A5").Value = "field"
$ws.Range("C5").Value = "string"
$ws.Range("A6").Value = "field"

# Bold header row (row 1) matches the style used on the other sheets
$ws.Range("A1:H1").Font.Bold = $true

# Column widths to match the authored layout
$ws.Columns.Item(1).ColumnWidth = 11.84375
$ws.Columns.Item(2).ColumnWidth = 16.4609375
$ws.Columns.Item(3).ColumnWidth = 8.84375
$ws.Columns.Item(4).ColumnWidth = 13.23046875
$ws.Columns.Item(5).ColumnWidth = 11.921875
$ws.Columns.Item(6).ColumnWidth = 11.3046875
$ws.Columns.Item(7).ColumnWidth = 12.53515625
$ws.Columns.Item(8).ColumnWidth = 15.69140625

# Page setup - portrait orientation like the other data sheets
$ws.PageSetup.Orientation = 1

# Selection/active cell as recorded in the source file
$ws.Range("C8").Select()
